$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time range values in column C
$ws.Range("C2").Value = "8:35-8:40"
$ws.Range("C3").Value = "8:40-8:45"

# Move the active selection to C14
$ws.Range("C14").Select()
